# Auto-generated edit script: refresh market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 55
$ws.Range("H55").Value = 310.2143
$ws.Range("I55").Value = 249.11111
$ws.Range("J55").Value = 420.2
$ws.Range("K55").Value = 249.11111
$ws.Range("L55").Value = 420.2
$ws.Range("M55").Value = -35.11111
$ws.Range("N55").Value = -848.2
# Row 112
$ws.Range("H112").Value = 2079.0625
$ws.Range("J112").Value = 2110.6453
$ws.Range("L112").Value = 6331.9359
$ws.Range("N112").Value = -8547.9359
# Row 132
$ws.Range("H132").Value = 9528436
$ws.Range("I132").Value = 11496984
$ws.Range("J132").Value = 13786.667
$ws.Range("K132").Value = 34490952
$ws.Range("L132").Value = 41360.001
$ws.Range("M132").Value = -34488422
$ws.Range("N132").Value = -46420.001

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 14972.857
$ws.Range("I2").Value = 659.4
$ws.Range("J2").Value = 50756.5
$ws.Range("K2").Value = 659.4
$ws.Range("L2").Value = 50756.5
$ws.Range("M2").Value = -546.4
$ws.Range("N2").Value = -50982.5
# Row 13
$ws.Range("H13").Value = 16666667
$ws.Range("I13").Value = 16666667
$ws.Range("K13").Value = 16666667
$ws.Range("M13").Value = -16666523
# Row 61
$ws.Range("H61").Value = 1853.5
$ws.Range("I61").Value = 1200
$ws.Range("J61").Value = 2507
$ws.Range("K61").Value = 1200
$ws.Range("L61").Value = 2507
$ws.Range("M61").Value = -988
$ws.Range("N61").Value = -2931
# Row 114
$ws.Range("H114").Value = 21696
$ws.Range("J114").Value = 21696
$ws.Range("L114").Value = 21696
$ws.Range("N114").Value = -30374
# Row 116
$ws.Range("H116").Value = 14972.857
$ws.Range("I116").Value = 659.4
$ws.Range("J116").Value = 50756.5
$ws.Range("K116").Value = 659.4
$ws.Range("L116").Value = 50756.5
$ws.Range("M116").Value = 1634.6
$ws.Range("N116").Value = -55344.5
# Row 122
$ws.Range("H122").Value = 1045.1
$ws.Range("I122").Value = 1045.1
$ws.Range("K122").Value = 3135.3
$ws.Range("M122").Value = -685.2999999999997
# Row 132
$ws.Range("H132").Value = 3840.9092
$ws.Range("I132").Value = 4164.8
$ws.Range("J132").Value = 3571
$ws.Range("K132").Value = 12494.4
$ws.Range("L132").Value = 10713
$ws.Range("M132").Value = -9964.400000000001
$ws.Range("N132").Value = -15773
# Row 136
$ws.Range("H136").Value = 1853.5
$ws.Range("I136").Value = 1200
$ws.Range("J136").Value = 2507
$ws.Range("K136").Value = 3600
$ws.Range("L136").Value = 7521
$ws.Range("M136").Value = -1050
$ws.Range("N136").Value = -12621

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 14972.857
$ws.Range("I3").Value = 659.4
$ws.Range("J3").Value = 50756.5
$ws.Range("K3").Value = 659.4
$ws.Range("L3").Value = 50756.5
$ws.Range("M3").Value = -545.4
$ws.Range("N3").Value = -50984.5
# Row 105
$ws.Range("H105").Value = 142859680
$ws.Range("I105").Value = 166669170
$ws.Range("J105").Value = 2800
$ws.Range("K105").Value = 166669170
$ws.Range("L105").Value = 2800
$ws.Range("M105").Value = -166667423
$ws.Range("N105").Value = -6294
# Row 134
$ws.Range("H134").Value = 8703.777
$ws.Range("I134").Value = 1257.8334
$ws.Range("J134").Value = 12426.75
$ws.Range("K134").Value = 3773.5002
$ws.Range("L134").Value = 37280.25
$ws.Range("M134").Value = -1238.5002
$ws.Range("N134").Value = -42350.25

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 1151.2354
$ws.Range("I58").Value = 855.12
$ws.Range("J58").Value = 1973.7778
$ws.Range("K58").Value = 855.12
$ws.Range("L58").Value = 1973.7778
$ws.Range("M58").Value = -652.12
$ws.Range("N58").Value = -2379.7778
# Row 132
$ws.Range("H132").Value = 4918.3447
$ws.Range("I132").Value = 5987.45
$ws.Range("J132").Value = 2542.5557
$ws.Range("K132").Value = 17962.35
$ws.Range("L132").Value = 7627.6671
$ws.Range("M132").Value = -15432.35
$ws.Range("N132").Value = -12687.6671
# Row 134
$ws.Range("H134").Value = 2678
$ws.Range("I134").Value = 3356.6667
$ws.Range("K134").Value = 10070.0001
$ws.Range("M134").Value = -7535.000100000001
# Row 136
$ws.Range("H136").Value = 1151.2354
$ws.Range("I136").Value = 855.12
$ws.Range("J136").Value = 1973.7778
$ws.Range("K136").Value = 2565.36
$ws.Range("L136").Value = 5921.3334
$ws.Range("M136").Value = -15.36000000000013
$ws.Range("N136").Value = -11021.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 74
$ws.Range("J12").Value = 66
$ws.Range("L12").Value = 198
$ws.Range("N12").Value = -544
# Row 22
$ws.Range("H22").Value = 2500
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 7500
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -7838
# Row 27
$ws.Range("H27").Value = 2500
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 7500
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -7704
# Row 55
$ws.Range("H55").Value = 3222.2222
$ws.Range("J55").Value = 3222.2222
$ws.Range("L55").Value = 9666.6666
$ws.Range("N55").Value = -10020.6666
# Row 107
$ws.Range("H107").Value = 748.3333
$ws.Range("I107").Value = 295
$ws.Range("J107").Value = 975
$ws.Range("K107").Value = 885
$ws.Range("L107").Value = 2925
$ws.Range("M107").Value = 1035
$ws.Range("N107").Value = -6765
# Row 131
$ws.Range("H131").Value = 13334444
$ws.Range("J131").Value = 1220.277
$ws.Range("L131").Value = 3660.831
$ws.Range("N131").Value = -13740.831

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Range("H126").Value = 2680.3333
$ws.Range("I126").Value = 3105.75
$ws.Range("J126").Value = 2340
$ws.Range("K126").Value = 9317.25
$ws.Range("L126").Value = 7020
$ws.Range("M126").Value = -6847.25
$ws.Range("N126").Value = -11960
# Row 132
$ws.Range("H132").Value = 3166.3845
$ws.Range("I132").Value = 2983.125
$ws.Range("J132").Value = 3459.6
$ws.Range("K132").Value = 8949.375
$ws.Range("L132").Value = 10378.8
$ws.Range("M132").Value = -6419.375
$ws.Range("N132").Value = -15438.8

$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Range("H3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").ClearContents()
# Row 5
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
# Row 7
$ws.Range("H7").Value = 3102.125
$ws.Range("I7").Value = 2900
$ws.Range("J7").Value = 3439
$ws.Range("K7").Value = 2900
$ws.Range("L7").Value = 3439
$ws.Range("M7").Value = -2788
$ws.Range("N7").Value = -3663
# Row 15
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
# Row 126
$ws.Range("H126").Value = 3102.125
$ws.Range("I126").Value = 2900
$ws.Range("J126").Value = 3439
$ws.Range("K126").Value = 8700
$ws.Range("L126").Value = 10317
$ws.Range("M126").Value = -6230
$ws.Range("N126").Value = -15257

$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 3496
$ws.Range("J26").Value = 498
$ws.Range("L26").Value = 498
$ws.Range("N26").Value = -1084
# Row 122
$ws.Range("H122").Value = 8967584
$ws.Range("I122").Value = 10402094
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 31206282
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -31203832
$ws.Range("N122").Value = -10600
# Row 136
$ws.Range("H136").Value = 574
$ws.Range("I136").Value = 559.9375
$ws.Range("J136").Value = 799
$ws.Range("K136").Value = 1679.8125
$ws.Range("L136").Value = 2397
$ws.Range("M136").Value = 870.1875
$ws.Range("N136").Value = -7497
# Row 137
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
